# Update the build timestamp embedded in the "build_version" strings
# throughout the workbook, from "17.29.55" to "18.05.36" (February 03 2026, EST).

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet: A2 (version banner) and A6 (recommended citation) ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = $wsAbout.Range("A2").Value().Replace($oldStamp, $newStamp)
$wsAbout.Range("A6").Value = $wsAbout.Range("A6").Value().Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet: build_version column (S), rows 2-10 ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 19).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = $val.Replace($oldStamp, $newStamp)
    }
}
